# Append new log rows (29-36) to the "Data" sheet, mirroring the
# ESTIM simulation log entries recorded on 2018.08.23 and 2018.08.24.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper-less approach: set date cells (column A) as text using NumberFormat "@"
# so Excel does not auto-convert the "yyyy.mm.dd"-looking text into a date
# serial number, then reset the style back to Normal (style index 0) so no
# extra cell style is introduced.

function Set-DateText($cellAddr, $text) {
    $rng = $ws.Range($cellAddr)
    $rng.NumberFormat = "@"
    $rng.Value = $text
    $rng.Style = "Normal"
}

$rows = @(
    @{ r = 29; A = "2018.08.23"; B = "19:03:56"; C = "RS"; D = 10; E = 250; F = 0.1;   G = 0.97;               H = 2975; I = 0.42; J = 0; K = "N/A"; L = "N/A" },
    @{ r = 30; A = "2018.08.23"; B = "19:17:27"; C = "RS"; D = 10; E = 250; F = 0.1;   G = 0.96;               H = 2975; I = 0.43; J = 0; K = "N/A"; L = "N/A" },
    @{ r = 31; A = "2018.08.24"; B = "15:00:27"; C = "RS"; D = 10; E = 250; F = 0.1;   G = 0.96;               H = 2975; I = 0.48; J = 0; K = "N/A"; L = "N/A" },
    @{ r = 32; A = "2018.08.24"; B = "15:09:13"; C = "RS"; D = 10; E = 250; F = 0.117; G = 0.91;               H = 2965; I = 0.44; J = 0; K = "N/A"; L = "N/A" },
    @{ r = 33; A = "2018.08.24"; B = "16:34:48"; C = "RS"; D = 10; E = 250; F = 0.1;   G = 0.9500000000000001; H = 2975; I = 0.44; J = 0; K = "N/A"; L = "N/A" },
    @{ r = 34; A = "2018.08.24"; B = "16:38:02"; C = "RS"; D = 10; E = 250; F = 0.1;   G = 0.9500000000000001; H = 2975; I = 0.43; J = 0; K = "N/A"; L = "N/A" },
    @{ r = 35; A = "2018.08.24"; B = "17:02:16"; C = "RS"; D = 10; E = 250; F = 0.1;   G = 0.9500000000000001; H = 2975; I = 0.45; J = 0; K = "N/A"; L = "N/A" },
    @{ r = 36; A = "2018.08.24"; B = "17:06:48"; C = "RS"; D = 10; E = 250; F = 0.1;   G = 0.97;               H = 2975; I = 0.41; J = 0; K = "N/A"; L = "N/A" }
)

foreach ($row in $rows) {
    $r = $row.r

    Set-DateText "A$r" $row.A

    $ws.Range("B$r").Value = $row.B
    $ws.Range("C$r").Value = $row.C
    $ws.Range("D$r").Value2 = $row.D
    $ws.Range("E$r").Value2 = $row.E
    $ws.Range("F$r").Value2 = $row.F
    $ws.Range("G$r").Value2 = $row.G
    $ws.Range("H$r").Value2 = $row.H
    $ws.Range("I$r").Value2 = $row.I
    $ws.Range("J$r").Value2 = $row.J
    $ws.Range("K$r").Value = $row.K
    $ws.Range("L$r").Value = $row.L
}
